$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.29347
$ws.Range("H2").Value = 0.8804099999999999
$ws.Range("I2").Value = 0.1501202107524681
$ws.Range("J2").Value = 0.1501202107524681
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.411745
$ws.Range("N2").Value = 16.235235
$ws.Range("O2").Value = 0.2701007085902594
$ws.Range("P2").Value = 0.2701007085902594
$ws.Range("Q2").Value = 1.58818480515
$ws.Range("R2").Value = 14.29366324635
$ws.Range("S2").Value = 0.04054757529796073
$ws.Range("T2").Value = 0.04054757529796073

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.29347
$ws.Range("H3").Value = 0.8804099999999999
$ws.Range("I3").Value = 0.1501202107524681
$ws.Range("J3").Value = 0.1501202107524681
$ws.Range("M3").Value = 4.518509
$ws.Range("O3").Value = 0.2255192147212155
$ws.Range("P3").Value = 0.2255192147212155
$ws.Range("Q3").Value = 1.32604683623
$ws.Range("R3").Value = 11.93442152607
$ws.Range("S3").Value = 0.03385499204267999
$ws.Range("T3").Value = 0.03385499204267998

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.29347
$ws.Range("H4").Value = 0.8804099999999999
$ws.Range("I4").Value = 0.1501202107524681
$ws.Range("J4").Value = 0.1501202107524681
$ws.Range("M4").Value = 10.10577266666666
$ws.Range("N4").Value = 30.317318
$ws.Range("O4").Value = 0.504380076688525
$ws.Range("P4").Value = 0.504380076688525
$ws.Range("Q4").Value = 2.965741104486666
$ws.Range("R4").Value = 26.69166994038
$ws.Range("S4").Value = 0.07571764341182742
$ws.Range("T4").Value = 0.07571764341182742

# Row 5
$ws.Range("I5").Value = 0.2180391153852712
$ws.Range("J5").Value = 0.2180391153852712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.411745
$ws.Range("N5").Value = 16.235235
$ws.Range("O5").Value = 0.2701007085902594
$ws.Range("P5").Value = 0.2701007085902594
$ws.Range("Q5").Value = 2.30672744361
$ws.Range("R5").Value = 20.76054699249
$ws.Range("S5").Value = 0.05889251956595509
$ws.Range("T5").Value = 0.05889251956595509

# Row 6
$ws.Range("I6").Value = 0.2180391153852712
$ws.Range("J6").Value = 0.2180391153852712
$ws.Range("M6").Value = 4.518509
$ws.Range("O6").Value = 0.2255192147212155
$ws.Range("P6").Value = 0.2255192147212155
$ws.Range("Q6").Value = 1.925990362535333
$ws.Range("S6").Value = 0.04917201008019487
$ws.Range("T6").Value = 0.04917201008019486

# Row 7
$ws.Range("I7").Value = 0.2180391153852712
$ws.Range("J7").Value = 0.2180391153852712
$ws.Range("M7").Value = 10.10577266666666
$ws.Range("N7").Value = 30.317318
$ws.Range("O7").Value = 0.504380076688525
$ws.Range("P7").Value = 0.504380076688525
$ws.Range("Q7").Value = 4.307531701712444
$ws.Range("R7").Value = 38.767785315412
$ws.Range("S7").Value = 0.1099745857391213
$ws.Range("T7").Value = 0.1099745857391213

# Row 8
$ws.Range("G8").Value = 1.235185333333333
$ws.Range("H8").Value = 3.705556
$ws.Range("I8").Value = 0.6318406738622607
$ws.Range("J8").Value = 0.6318406738622606
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.411745
$ws.Range("N8").Value = 16.235235
$ws.Range("O8").Value = 0.2701007085902594
$ws.Range("P8").Value = 0.2701007085902594
$ws.Range("Q8").Value = 6.68450805174
$ws.Range("R8").Value = 60.16057246566
$ws.Range("S8").Value = 0.1706606137263436
$ws.Range("T8").Value = 0.1706606137263436

# Row 9
$ws.Range("G9").Value = 1.235185333333333
$ws.Range("H9").Value = 3.705556
$ws.Range("I9").Value = 0.6318406738622607
$ws.Range("J9").Value = 0.6318406738622606
$ws.Range("M9").Value = 4.518509
$ws.Range("O9").Value = 0.2255192147212155
$ws.Range("P9").Value = 0.2255192147212155
$ws.Range("Q9").Value = 5.581196045334667
$ws.Range("R9").Value = 50.230764408012
$ws.Range("S9").Value = 0.1424922125983407
$ws.Range("T9").Value = 0.1424922125983406

# Row 10
$ws.Range("G10").Value = 1.235185333333333
$ws.Range("H10").Value = 3.705556
$ws.Range("I10").Value = 0.6318406738622607
$ws.Range("J10").Value = 0.6318406738622606
$ws.Range("M10").Value = 10.10577266666666
$ws.Range("N10").Value = 30.317318
$ws.Range("O10").Value = 0.504380076688525
$ws.Range("P10").Value = 0.504380076688525
$ws.Range("Q10").Value = 12.48250217986755
$ws.Range("R10").Value = 112.342519618808
$ws.Range("S10").Value = 0.3186878475375764
$ws.Range("T10").Value = 0.3186878475375763
